$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1616541353383459
$ws.Range("C2").Value = 0.6240601503759399
$ws.Range("J2").Value = 0.003759398496240601
$ws.Range("P2").Value = 0.1203007518796992
$ws.Range("S2").Value = 0.09022556390977443
$ws.Range("J3").Value = 0.02409638554216868
$ws.Range("P3").Value = 0.7108433734939759
$ws.Range("S3").Value = 0.2650602409638554
$ws.Range("J4").Value = 0.0303030303030303
$ws.Range("P4").Value = 0.6363636363636364
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.08426966292134831
$ws.Range("D6").Value = 0.01123595505617977
$ws.Range("F6").Value = 0.03370786516853932
$ws.Range("J6").Value = 0.2191011235955056
$ws.Range("O6").Value = 0.02247191011235955
$ws.Range("Q6").Value = 0.1685393258426966
$ws.Range("R6").Value = 0.09550561797752809
$ws.Range("S6").Value = 0.3651685393258427
$ws.Range("B7").Value = 0.06741573033707865
$ws.Range("D7").Value = 0.01685393258426966
$ws.Range("F7").Value = 0.03932584269662921
$ws.Range("J7").Value = 0.1348314606741573
$ws.Range("O7").Value = 0.02247191011235955
$ws.Range("Q7").Value = 0.2078651685393259
$ws.Range("R7").Value = 0.05617977528089887
$ws.Range("S7").Value = 0.4550561797752809
$ws.Range("B8").Value = 0.09921671018276762
$ws.Range("D8").Value = 0.02610966057441253
$ws.Range("F8").Value = 0.04960835509138381
$ws.Range("J8").Value = 0.1174934725848564
$ws.Range("O8").Value = 0.007832898172323759
$ws.Range("Q8").Value = 0.1827676240208877
$ws.Range("R8").Value = 0.1148825065274151
$ws.Range("S8").Value = 0.402088772845953
$ws.Range("B9").Value = 0.07653061224489796
$ws.Range("D9").Value = 0.01020408163265306
$ws.Range("F9").Value = 0.0663265306122449
$ws.Range("J9").Value = 0.1275510204081633
$ws.Range("O9").Value = 0.01530612244897959
$ws.Range("Q9").Value = 0.1683673469387755
$ws.Range("R9").Value = 0.06122448979591837
$ws.Range("S9").Value = 0.4744897959183674
$ws.Range("B10").Value = 0.1130856219709208
$ws.Range("D10").Value = 0.01373182552504039
$ws.Range("E10").Value = 0.002423263327948304
$ws.Range("F10").Value = 0.06946688206785137
$ws.Range("J10").Value = 0.1017770597738288
$ws.Range("O10").Value = 0.01373182552504039
$ws.Range("Q10").Value = 0.2237479806138934
$ws.Range("R10").Value = 0.08966074313408724
$ws.Range("S10").Value = 0.3723747980613893
$ws.Range("G11").Value = 0.156794425087108
$ws.Range("J11").Value = 0.10801393728223
$ws.Range("K11").Value = 0.2229965156794425
$ws.Range("L11").Value = 0.5052264808362369
$ws.Range("S11").Value = 0.006968641114982578
$ws.Range("G12").Value = 0.7218543046357616
$ws.Range("J12").Value = 0.2185430463576159
$ws.Range("K12").Value = 0.01324503311258278
$ws.Range("L12").Value = 0.02649006622516556
$ws.Range("S12").Value = 0.01986754966887417
$ws.Range("G13").Value = 0.6944444444444444
$ws.Range("J13").Value = 0.3055555555555556
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("J14").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.01123595505617977
$ws.Range("H15").Value = 0.151685393258427
$ws.Range("I15").Value = 0.07865168539325842
$ws.Range("J15").Value = 0.3820224719101123
$ws.Range("K15").Value = 0.02808988764044944
$ws.Range("M15").Value = 0.02247191011235955
$ws.Range("O15").Value = 0.0449438202247191
$ws.Range("S15").Value = 0.2808988764044944
$ws.Range("H16").Value = 0.1676646706586826
$ws.Range("I16").Value = 0.09580838323353294
$ws.Range("J16").Value = 0.4011976047904192
$ws.Range("K16").Value = 0.1317365269461078
$ws.Range("M16").Value = 0.005988023952095809
$ws.Range("O16").Value = 0.02994011976047904
$ws.Range("S16").Value = 0.1676646706586826
$ws.Range("F17").Value = 0.01789709172259508
$ws.Range("H17").Value = 0.1543624161073825
$ws.Range("I17").Value = 0.1140939597315436
$ws.Range("J17").Value = 0.4519015659955257
$ws.Range("K17").Value = 0.07606263982102908
$ws.Range("M17").Value = 0.01565995525727069
$ws.Range("N17").Value = 0.002237136465324385
$ws.Range("O17").Value = 0.0447427293064877
$ws.Range("S17").Value = 0.1230425055928412
$ws.Range("F18").Value = 0.02072538860103627
$ws.Range("H18").Value = 0.1917098445595855
$ws.Range("I18").Value = 0.08808290155440414
$ws.Range("J18").Value = 0.4404145077720207
$ws.Range("K18").Value = 0.08290155440414508
$ws.Range("M18").Value = 0.0155440414507772
$ws.Range("O18").Value = 0.04663212435233161
$ws.Range("S18").Value = 0.1139896373056995
$ws.Range("F19").Value = 0.004838709677419355
$ws.Range("H19").Value = 0.1862903225806452
$ws.Range("I19").Value = 0.08064516129032258
$ws.Range("J19").Value = 0.3959677419354839
$ws.Range("K19").Value = 0.1145161290322581
$ws.Range("M19").Value = 0.01774193548387097
$ws.Range("N19").Value = 0.001612903225806452
$ws.Range("O19").Value = 0.06854838709677419
$ws.Range("S19").Value = 0.1298387096774194
